$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must remain text even when it looks like a number.
# Excel auto-converts plain numeric-looking strings to numbers on assignment,
# so we prefix with an apostrophe (forces text entry) and then restore the
# cell's original style so no formatting side effects are introduced.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

# Row 2
$ws.Range("D2").Value = "63.657.36"
$ws.Range("E2").Value = "  -1.80%  "

# Row 3
$ws.Range("D3").Value = "3.125.33"
$ws.Range("E3").Value = "  -1.28%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
Set-TextValue $ws.Range("D5") "600.73"
$ws.Range("E5").Value = "  -2.27%  "

# Row 6
Set-TextValue $ws.Range("D6") "142.25"
$ws.Range("E6").Value = "  -4.64%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").Value = "3.118.47"
$ws.Range("E8").Value = "  -1.61%  "

# Row 9
$ws.Range("E9").Value = "  -1.00%  "

# Row 10
$ws.Range("E10").Value = "  -2.75%  "

# Row 11
Set-TextValue $ws.Range("D11") "5.35"
$ws.Range("E11").Value = "  -3.16%  "

# Row 12
$ws.Range("E12").Value = "  -2.07%  "

# Row 13
$ws.Range("E13").Value = "  -2.81%  "

# Row 14
Set-TextValue $ws.Range("D14") "34.94"
$ws.Range("E14").Value = "  -2.48%  "

# Row 15
$ws.Range("D15").Value = "3.638.84"
$ws.Range("E15").Value = "  -1.24%  "

# Row 17
$ws.Range("D17").Value = "63.748.46"
$ws.Range("E17").Value = "  -1.52%  "

# Row 18
$ws.Range("D18").Value = "3.116.87"
$ws.Range("E18").Value = "  -1.36%  "

# Row 19
Set-TextValue $ws.Range("D19") "6.80"
$ws.Range("E19").Value = "  -2.12%  "

# Row 20
Set-TextValue $ws.Range("D20") "483.23"
$ws.Range("E20").Value = "  -0.12%  "

# Row 21
$ws.Range("E21").Value = "  -0.82%  "

# Row 22
$ws.Range("E22").Value = "  -1.68%  "

# Row 23
Set-TextValue $ws.Range("D23") "7.62"
$ws.Range("E23").Value = "  -5.74%  "

# Row 24
Set-TextValue $ws.Range("D24") "86.61"
$ws.Range("E24").Value = "  +2.76%  "

# Row 25
Set-TextValue $ws.Range("D25") "13.40"
$ws.Range("E25").Value = "  -2.71%  "

# Row 26
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("E27").Value = "  -4.01%  "

# Row 28
Set-TextValue $ws.Range("D28") "8.23"
$ws.Range("E28").Value = "  -4.81%  "

# Row 29
Set-TextValue $ws.Range("D29") "7.02"
$ws.Range("E29").Value = "  -1.60%  "

# Row 30
$ws.Range("E30").Value = "  -2.95%  "

# Row 31
$ws.Range("E31").Value = "  +1.86%  "

# Row 32
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D32") "1.00"
$ws.Range("E32").Value = "  -0.01%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D33") "0.111"
$ws.Range("E33").Value = "  -8.29%  "

# Row 34
Set-TextValue $ws.Range("D34") "2.64"
$ws.Range("E34").Value = "  -3.87%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.11"
$ws.Range("E35").Value = "  -3.31%  "

# Row 36
$ws.Range("E36").Value = "  -1.55%  "

# Row 37
Set-TextValue $ws.Range("D37") "52.54"
$ws.Range("E37").Value = "  -1.42%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0742"
$ws.Range("E38").Value = "  -7.02%  "

# Row 39
Set-TextValue $ws.Range("D39") "2.95"
$ws.Range("E39").Value = "  -10.76%  "

# Row 40
Set-TextValue $ws.Range("D40") "436.42"
$ws.Range("E40").Value = "  -5.93%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.0393"
$ws.Range("E41").Value = "  -2.47%  "

# Row 42
$ws.Range("E42").Value = "  -0.93%  "

# Row 43
$ws.Range("E43").Value = "  -2.24%  "

# Row 44
$ws.Range("D44").Value = "2.870.74"
$ws.Range("E44").Value = "  +0.46%  "

# Row 45
$ws.Range("E45").Value = "  -4.50%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.19"
$ws.Range("E46").Value = "  -6.64%  "

# Row 47
$ws.Range("E47").Value = "  +0.05%  "

# Row 48
$ws.Range("E48").Value = "  -5.07%  "

# Row 49
Set-TextValue $ws.Range("D49") "25.76"
$ws.Range("E49").Value = "  -3.69%  "

# Row 50
$ws.Range("E50").Value = "  -1.15%  "

# Row 51
Set-TextValue $ws.Range("D51") "121.47"
$ws.Range("E51").Value = "  +1.01%  "

